$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New delivery: widen the first four columns (Clasificacion, Tipo de
# producto o servicio, Nombre Producto/Servicio, Nombre Etiqueta) so the
# newly imported data is readable.
# NOTE: the host's ColumnWidth setter quantizes to 1/6-character steps
# (stored width = round(chars*6)/6 + 5/6), so we back-solve the input that
# lands closest to each target stored width.
$ws.Columns.Item(1).ColumnWidth = 16.330729166666668   # -> stored width 17.1640625 (A)
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668   # -> stored width 23         (B)
$ws.Columns.Item(3).ColumnWidth = 33.666666666666664   # -> stored width 34.5       (C)
$ws.Columns.Item(4).ColumnWidth = 18.830729166666668   # -> stored width 19.6640625 (D)

# Move the active selection to H7
$ws.Range("H7").Select()
